# Gender-detection workbook: add a "Precission and Recall" sheet summarizing
# accuracy / recall / precision from a confusion matrix.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 (it becomes the active / selected tab).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws2.Name = "Precission and Recall"

# --- text labels, entered in the order that reproduces the original shared-string table ---
$ws2.Range("A3").Value = "Recall"
$ws2.Range("A4").Value = "Precision"
$ws2.Range("A2").Value = "Accuracy"
$ws2.Range("D3").Value = " "
$ws2.Range("E3").Value = "Negative cases"
$ws2.Range("E4").Value = "Positive cases"
$ws2.Range("F2").Value = "Predicted Negative"
$ws2.Range("G2").Value = "Predicted Postiive"
$ws2.Range("A15").Value = "TOTAL"
$ws2.Range("A8").Value = "True positive 294"
$ws2.Range("A9").Value = "True negative 194"
$ws2.Range("A10").Value = "False positive 100"
$ws2.Range("A11").Value = "False negative 412"

# --- confusion-matrix raw counts ---
$ws2.Range("B8").Value = 294
$ws2.Range("B9").Value = 194
$ws2.Range("B10").Value = 100
$ws2.Range("B11").Value = 412
$ws2.Range("B15").Formula = "=SUM(B8:B11)"

# --- confusion-matrix grid (predicted negative / predicted positive) ---
$ws2.Range("F3").Formula = "=B9"
$ws2.Range("G3").Formula = "=B10"
$ws2.Range("F4").Formula = "=B11"
$ws2.Range("G4").Formula = "=B8"

# --- accuracy / recall / precision formulas ---
$ws2.Range("B2").Formula = "=(F3+G4)/B15"
$ws2.Range("B3").Formula = "=G4/(G4+F4)"
$ws2.Range("B4").Formula = "=G4/(G4+G3)"
$ws2.Range("B2:B4").NumberFormat = "0.00%"

# --- column widths (target stored widths: A=23.6640625, E=19.33203125, F=18.1640625, G=18) ---
$ws2.Columns.Item(1).ColumnWidth = 23.6640625 - 0.8333333333333334
$ws2.Columns.Item(5).ColumnWidth = 19.33203125 - 0.8333333333333334
$ws2.Columns.Item(6).ColumnWidth = 18.1640625 - 0.8333333333333334
$ws2.Columns.Item(7).ColumnWidth = 18 - 0.8333333333333334

# --- page setup, matching Sheet1 ---
$ps2 = $ws2.PageSetup
$ps2.LeftMargin = 54
$ps2.RightMargin = 54
$ps2.TopMargin = 72
$ps2.BottomMargin = 72
$ps2.HeaderMargin = 36
$ps2.FooterMargin = 36
$ps2.PaperSize = 9
$ps2.Orientation = 1

# --- selection / active view state ---
$ws2.Activate()
$ws2.Range("E31").Select()
